# Auto-generated edit script applying the cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Sat Dec 30 08:23:38 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.889.42"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").Value = "2.289.24"
$ws.Range("E3").Value = "  -2.30%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'317.03"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'102.82"
$ws.Range("E6").Value = "  -4.27%  "
$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.604"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("D10").Value = "'38.98"
$ws.Range("E10").Value = "  -4.71%  "
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").Value = "'8.31"
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").Value = "'0.966"
$ws.Range("E14").Value = "  -2.49%  "
$ws.Range("D15").Value = "'15.29"
$ws.Range("E15").Value = "  -4.08%  "
$ws.Range("D16").Value = "2.636.39"
$ws.Range("E16").Value = "  -2.36%  "
$ws.Range("D17").Value = "2.301.86"
$ws.Range("E17").Value = "  -3.75%  "
$ws.Range("D18").Value = "41.985.02"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").Value = "'7.58"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").Value = "'284.81"
$ws.Range("E21").Value = "  +11.22%  "
$ws.Range("D22").Value = "'73.78"
$ws.Range("E22").Value = "  -2.97%  "
$ws.Range("D23").Value = "'3.57"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("D24").Value = "'2.26"
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("D25").Value = "'9.93"
$ws.Range("E25").Value = "  +6.07%  "
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("D27").Value = "'10.78"
$ws.Range("E27").Value = "  -5.03%  "
$ws.Range("D28").Value = "'23.18"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").Value = "'2.27"
$ws.Range("E29").Value = "  +3.13%  "
$ws.Range("D30").Value = "'163.26"
$ws.Range("E30").Value = "  -6.17%  "
$ws.Range("D31").Value = "'34.88"
$ws.Range("E31").Value = "  -4.96%  "
$ws.Range("D32").Value = "'0.0879"
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("E33").Value = "  +2.36%  "
$ws.Range("D34").Value = "'5.84"
$ws.Range("E34").Value = "  -3.77%  "
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("E36").Value = "  -7.68%  "
$ws.Range("D37").Value = "'4.58"
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.88"
$ws.Range("E38").Value = "  +8.89%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0348"
$ws.Range("E39").Value = "  -3.70%  "
$ws.Range("D40").Value = "'3.60"
$ws.Range("E40").Value = "  -7.28%  "
$ws.Range("D41").Value = "'101.91"
$ws.Range("E41").Value = "  +18.67%  "
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("D43").Value = "'70.02"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.225"
$ws.Range("E45").Value = "  -4.45%  "
$ws.Range("D46").Value = "'115.77"
$ws.Range("E46").Value = "  +4.10%  "
$ws.Range("D47").Value = "'11.95"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'9.06"
$ws.Range("E48").Value = "  -1.00%  "
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "'76.59"
$ws.Range("E49").Value = "  +2.58%  "
$ws.Range("D50").Value = "'5.31"
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("D51").Value = "'1.27"
$ws.Range("E51").Value = "  -0.51%  "
